$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Gulf_of_Oman: drop the saved scroll position / old selection, select G10,
# and resize the first two columns.
# ---------------------------------------------------------------------------
$wsGulf = $wb.Worksheets.Item("Gulf_of_Oman")
$wsGulf.Range("G10").Select()
$wsGulf.Columns.Item(1).ColumnWidth = 7
$wsGulf.Columns.Item(2).ColumnWidth = 7.833333333333332

# ---------------------------------------------------------------------------
# F1: select column B (whole column) instead of the full-sheet selection.
# ---------------------------------------------------------------------------
$wsF1 = $wb.Worksheets.Item("F1")
$wsF1.Columns("B:B").Select()

# ---------------------------------------------------------------------------
# F3: select column A (whole column) and size the first two columns.
# ---------------------------------------------------------------------------
$wsF3 = $wb.Worksheets.Item("F3")
$wsF3.Columns("A:A").Select()
$wsF3.Columns.Item(1).ColumnWidth = 8.5
$wsF3.Columns.Item(2).ColumnWidth = 10.333333333333332

# ---------------------------------------------------------------------------
# F2: relabel the header row (Age / % dolo wt.), replace the data column
# with the dolomite wt% series, widen column B, and finish by selecting
# columns A:B so F2 ends up the active sheet/tab.
# ---------------------------------------------------------------------------
$wsF2 = $wb.Worksheets.Item("F2")

$wsF2.Cells.Item(1, 1).Value = "Age"

$wsF2.Cells.Item(1, 2).ClearFormats()
$wsF2.Cells.Item(1, 2).Value = "% dolo wt."

# give the shared bold header font back its family attribute
$wsF2.Cells.Item(1, 1).Font.Family = 2

$newVals = @(1.88205185,1.338571059,2.069281397,1.722448653,1.882497636,1.516210175,2.320288956,2.345055487,2.195841451,1.623097045,1.945657524,2.809948052,1.341770282,1.949810709,1.371766565,1.650652135,2.046868834,2.551221989,2.570715534,2.509357539,1.692957869,2.370255905,1.567587843,1.548859665,1.760981239,2.081061242,2.259838795,1.46374012,1.511509747,2.087132522,2.678952966,1.466312984,2.788160422,2.735410465,2.981581366,1.445819735,2.540970948,2.157169555,2.95033958,2.19041067,1.820700687,1.660867393,2.575671375,2.459063884,2.356243504,2.006101746,1.026688165,1.629269239,1.335483111,2.920183684,1.291382976,2.912508215,2.542421145,1.505551693,2.845280981,2.16069651,2.389052518,1.491426207,1.93656674,1.527163937)

for ($i = 0; $i -lt $newVals.Length; $i++) {
    $wsF2.Cells.Item($i + 2, 2).Value = $newVals[$i]
}

$wsF2.Columns.Item(2).ColumnWidth = 9.166666666666666

$wsF2.Columns("A:B").Select()
